# Applies:
#  1) Adds <w:ind w:firstLine="708"/> (i.e. a 0.5"/708-twip first-line indent)
#     to the paragraph properties of the "justified" body paragraphs
#     (those with jc=both), except for the sub-heading-styled one
#     ("Verificando o histórico de alterações").
#  2) Moves the (hidden) "_GoBack" bookmark from the very end of the
#     document (after the last image) to the start of the
#     "Lembrando que..." paragraph.

$d = $word.ActiveDocument

# 1-based paragraph indexes (in document order) of the body paragraphs
# that receive the new first-line indent.
$targets = @(2, 4, 8, 12, 15, 21, 24, 26)

# 708 twips = 35.4 points (Word's object model works in points; 1 pt = 20 twips)
$firstLinePoints = 708 / 20

foreach ($t in $targets) {
    $p = $d.Paragraphs.Item($t)
    $p.Format.FirstLineIndent = $firstLinePoints
}

# Move the "_GoBack" bookmark to the start of paragraph 26
# ("Lembrando que como a partir do momento...").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$targetParagraph = $d.Paragraphs.Item(26)
$startPos = $targetParagraph.Range.Start
$bookmarkRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
